$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.468.57'
$ws.Range("E2").Value = '  -0.87%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.374.71'
$ws.Range("E3").Value = '  -0.87%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '507.87'
$ws.Range("E5").Value = '  +0.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.50'
$ws.Range("E6").Value = '  -1.91%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  -1.47%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.384.98'
$ws.Range("E9").Value = '  -1.04%  '

$ws.Range("E11").Value = '  +0.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.86'
$ws.Range("E12").Value = '  +5.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.330'
$ws.Range("E13").Value = '  +2.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.798.27'
$ws.Range("E14").Value = '  -0.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.437.39'
$ws.Range("E15").Value = '  -0.87%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.58'
$ws.Range("E16").Value = '  -1.34%  '

$ws.Range("E17").Value = '  -0.64%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.390.15'
$ws.Range("E18").Value = '  -0.24%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.05'
$ws.Range("E19").Value = '  -1.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.05'
$ws.Range("E20").Value = '  +0.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '309.57'
$ws.Range("E21").Value = '  -0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.30'
$ws.Range("E22").Value = '  -0.66%  '

$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.16'
$ws.Range("E24").Value = '  +1.49%  '

$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("E26").Value = '  -1.43%  '

$ws.Range("E27").Value = '  -2.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.25'
$ws.Range("E28").Value = '  -2.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.59'
$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("E30").Value = '  -1.44%  '

$ws.Range("E31").Value = '  -1.52%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.84'
$ws.Range("E32").Value = '  -1.59%  '

$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.08'
$ws.Range("E35").Value = '  -3.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.65'
$ws.Range("E36").Value = '  -1.68%  '

$ws.Range("E37").Value = '  -0.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.71'
$ws.Range("E38").Value = '  -3.29%  '

$ws.Range("E39").Value = '  +2.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.44'
$ws.Range("E40").Value = '  -0.81%  '

$ws.Range("E41").Value = '  -3.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.39'
$ws.Range("E42").Value = '  +1.13%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.96'
$ws.Range("E43").Value = '  -0.88%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '126.11'
$ws.Range("E44").Value = '  -4.91%  '

$ws.Range("E45").Value = '  -0.64%  '

$ws.Range("E46").Value = '  -1.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '239.94'
$ws.Range("E47").Value = '  -4.88%  '

$ws.Range("E48").Value = '  -1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0207'
$ws.Range("E49").Value = '  -1.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.03'
$ws.Range("E50").Value = '  -1.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.952'
$ws.Range("E51").Value = '  -0.14%  '
